$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3963402211666107
$ws.Range("B1").Value = 1.480807542800903
$ws.Range("C1").Value = 3.682619094848633
$ws.Range("D1").Value = 3.275521516799927
$ws.Range("E1").Value = 0.8217142224311829
